$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("StringLocalizations_BasicText")
$ws.Activate()

# The sheet currently has data rows through row 64 (A1:E64).
# We append 4 new rows (65-68), each holding a localization KEY (column A)
# and its EN-GB value (column B), with placeholder "XXXX" in columns C-E
# (matching the pattern used for every other row on this sheet).

# Copy the formatting (wrap-text style) of the last existing row down onto
# the new rows first, so the new cells pick up style index 1 just like the
# surrounding rows.
$ws.Range("A64:E64").Copy($ws.Range("A65:E68")) | Out-Null

# Row 65: INCIDENT_AVAILABLE / Available
$ws.Cells.Item(65, 1).Value2 = "INCIDENT_AVAILABLE"
$ws.Cells.Item(65, 2).Value2 = "Available"
$ws.Cells.Item(65, 3).Value2 = "XXXX"
$ws.Cells.Item(65, 4).Value2 = "XXXX"
$ws.Cells.Item(65, 5).Value2 = "XXXX"
# The key cell in this new row was typed without inheriting the wrap-text
# formatting used elsewhere, so strip its formatting back to default.
$ws.Cells.Item(65, 1).ClearFormats()

# Row 66: BASIC_TEXT_SATISFACTION / {0}% Satisfaction
$ws.Cells.Item(66, 1).Value2 = "BASIC_TEXT_SATISFACTION"
$ws.Cells.Item(66, 2).Value2 = "{0}% Satisfaction"
$ws.Cells.Item(66, 3).Value2 = "XXXX"
$ws.Cells.Item(66, 4).Value2 = "XXXX"
$ws.Cells.Item(66, 5).Value2 = "XXXX"

# Row 67: BASIC_TEXT_SATISFACTION_END_TURN / {0} Ignored Case{1}, -{2}% Satisfaction
$ws.Cells.Item(67, 1).Value2 = "BASIC_TEXT_SATISFACTION_END_TURN"
$ws.Cells.Item(67, 2).Value2 = "{0} Ignored Case{1}, -{2}% Satisfaction"
$ws.Cells.Item(67, 3).Value2 = "XXXX"
$ws.Cells.Item(67, 4).Value2 = "XXXX"
$ws.Cells.Item(67, 5).Value2 = "XXXX"

# Row 68: BASIC_TEXT_NO_IGNORED_CASES / 0 Ignored Cases, No Satisfaction Change
$ws.Cells.Item(68, 1).Value2 = "BASIC_TEXT_NO_IGNORED_CASES"
$ws.Cells.Item(68, 2).Value2 = "0 Ignored Cases, No Satisfaction Change"
$ws.Cells.Item(68, 3).Value2 = "XXXX"
$ws.Cells.Item(68, 4).Value2 = "XXXX"
$ws.Cells.Item(68, 5).Value2 = "XXXX"

# Move the visible selection/scroll position down to the newly added data,
# matching where the author ended up after adding these rows.
[void]$ws.Range("A68").Select()
